# Update the roster to match the new master-file de-duplicated list.
# Each data row (2-12) is updated in place with the new Name/Title/Office
# values (City stays "Vancouver, BC" throughout), and the old row 13 is
# removed entirely since the list shrank by one entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Name = "David Malkin";         Title = "Licensed Administrator"; Office = "RE/MAX Select Properties" },
    @{ Row = 3;  Name = "Shahina Najmidinova";   Title = "Licensed Administrator"; Office = "RE/MAX City Realty" },
    @{ Row = 4;  Name = "Leanne Yan";            Title = "Sales Representative";   Office = "RE/MAX Real Estate Services" },
    @{ Row = 5;  Name = "Nicola Campbell";       Title = "Sales Representative";   Office = "RE/MAX Crest Realty (South Granville)" },
    @{ Row = 6;  Name = "Hazem Sultan";          Title = "Real Estate Advisor";    Office = "RE/MAX Crest Realty (South Granville)" },
    @{ Row = 7;  Name = "Meet Dusange";          Title = "Sales Representative";   Office = "RE/MAX City Realty" },
    @{ Row = 8;  Name = "Khush Grewal";          Title = "Realtor";                Office = "RE/MAX Elevate" },
    @{ Row = 9;  Name = "Gabe Bandel";           Title = "Licensed Administrator"; Office = "RE/MAX Select Realty" },
    @{ Row = 10; Name = "James L Wang";          Title = "Sales Representative";   Office = "RE/MAX City Realty" },
    @{ Row = 11; Name = "Lina Rached";           Title = "Sales Representative";   Office = "RE/MAX Crest Realty (South Granville)" },
    @{ Row = 12; Name = "Rosalee McRae";         Title = "Broker";                 Office = "RE/MAX Select Properties" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Title
    $ws.Cells.Item($r, 4).Value = $entry.Office
}

# Remove the now-extra 13th row (list went from 12 entries to 11).
$ws.Rows.Item(13).Delete()
